$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.644.65'
$ws.Range("E2").Value = '  +0.91%  '

# Row 3
$ws.Range("D3").Value = '2.246.75'
$ws.Range("E3").Value = '  +0.10%  '

# Row 4
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '306.86'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '94.95'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.17%  '

# Row 7
$ws.Range("E7").Value = '  -0.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  -2.25%  '

# Row 10
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '34.99'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +0.41%  '

# Row 11
$ws.Range("E11").Value = '  -1.24%  '

# Row 12
$ws.Range("E12").Value = '  -0.04%  '

# Row 13
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +0.22%  '

# Row 14
$ws.Range("D14").Value = '2.591.16'
$ws.Range("E14").Value = '  +0.17%  '

# Row 15
$ws.Range("D15").Value = '2.238.76'
$ws.Range("E15").Value = '  -4.12%  '

# Row 16
$ws.Range("E16").Value = '  +0.05%  '

# Row 17
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '13.58'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +0.19%  '

# Row 18
$ws.Range("D18").Value = '44.430.97'
$ws.Range("E18").Value = '  +1.01%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0936'
$ws.Range("E19").Value = '  -3.03%  '

# Row 20
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '6.18'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -3.44%  '

# Row 21
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.73'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -3.57%  '

# Row 22
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '65.33'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '237.40'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.56%  '

# Row 24
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.96'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.01%  '

# Row 25
$ws.Range("E25").Value = '  -1.57%  '

# Row 26
$ws.Range("E26").Value = '  -0.17%  '

# Row 28
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '9.78'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -1.82%  '

# Row 29
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '37.02'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -3.63%  '

# Row 30
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '5.89'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +0.16%  '

# Row 31
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '19.99'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.52%  '

# Row 32
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '147.29'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -3.91%  '

# Row 33
$ws.Range("E33").Value = '  -1.64%  '

# Row 34
$ws.Range("E34").Value = '  +0.00%  '

# Row 35
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '3.20'
$ws.Range("D35").Style = 'Normal'

# Row 36
$ws.Range("E36").Value = '  +1.30%  '

# Row 37
$ws.Range("E37").Value = '  -1.64%  '

# Row 38
$ws.Range("E38").Value = '  +5.03%  '

# Row 39
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '15.19'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +5.15%  '

# Row 40
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '3.34'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -5.47%  '

# Row 41
$ws.Range("E41").Value = '  -1.62%  '

# Row 42
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.0300'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("E43").Value = '  +0.03%  '

# Row 44
$ws.Range("D44").Value = '1.810.83'
$ws.Range("E44").Value = '  +3.60%  '

# Row 45
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '1.77'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +12.24%  '

# Row 46
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '81.84'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.26%  '

# Row 47
$ws.Range("E47").Value = '  -2.03%  '

# Row 48
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '98.36'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.74%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '68.94'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +2.08%  '

# Row 50
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '4.81'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -2.44%  '

# Row 51
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '54.07'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -1.11%  '
